$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-detected as a number by Excel
# must have their number format forced to Text ("@") first, so the literal
# string (including e.g. trailing zeros) is preserved, matching the source data.
$textCells = @(
    'D5',
    'D6',
    'D8',
    'D11',
    'D15',
    'D16',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D30',
    'D32',
    'D33',
    'D34',
    'D37',
    'D38',
    'D40',
    'D41',
    'D42',
    'D43',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '98.167.86'
$ws.Range('E2').Value = '  +5.54%  '
$ws.Range('D3').Value = '3.140.26'
$ws.Range('E3').Value = '  +1.60%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '241.74'
$ws.Range('E5').Value = '  +2.64%  '
$ws.Range('D6').Value = '608.78'
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  -1.47%  '
$ws.Range('D8').Value = '0.382'
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '3.134.19'
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').Value = '0.787'
$ws.Range('E11').Value = '  -4.67%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '97.547.61'
$ws.Range('E13').Value = '  +4.97%  '
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('D15').Value = '5.42'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '33.76'
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('D17').Value = '3.719.19'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = '3.137.84'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '522.10'
$ws.Range('E19').Value = '  +18.74%  '
$ws.Range('D20').Value = '3.38'
$ws.Range('E20').Value = '  -7.78%  '
$ws.Range('D21').Value = '14.45'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').Value = '5.62'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('D23').Value = '0.0000190'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('D24').Value = '8.68'
$ws.Range('E24').Value = '  -3.66%  '
$ws.Range('D25').Value = '88.58'
$ws.Range('E25').Value = '  +3.65%  '
$ws.Range('D26').Value = '5.42'
$ws.Range('E26').Value = '  -4.14%  '
$ws.Range('D27').Value = '11.51'
$ws.Range('E27').Value = '  -8.81%  '
$ws.Range('D28').Value = '3.301.27'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = '0.238'
$ws.Range('E30').Value = '  -4.53%  '
$ws.Range('E31').Value = '  -3.31%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.121'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').Value = '8.87'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('E35').Value = '  +3.00%  '
$ws.Range('E36').Value = '  -5.30%  '
$ws.Range('D37').Value = '7.17'
$ws.Range('E37').Value = '  -9.14%  '
$ws.Range('D38').Value = '24.34'
$ws.Range('E38').Value = '  +1.55%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '0.432'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '465.03'
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('D42').Value = '1.20'
$ws.Range('E42').Value = '  -6.12%  '
$ws.Range('D43').Value = '3.50'
$ws.Range('E43').Value = '  -11.71%  '
$ws.Range('D45').Value = '3.07'
$ws.Range('E45').Value = '  -5.62%  '
$ws.Range('D46').Value = '162.55'
$ws.Range('E46').Value = '  +2.22%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '1.90'
$ws.Range('E47').Value = '  +2.94%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '0.691'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D49').Value = '4.47'
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('D50').Value = '44.17'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  +0.06%  '
